$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-18 16:43:29"
$wsZhCn.Range("H4").Value = "2016-03-18 16:44:15"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-18 16:43:37"
$wsDeDe.Range("H4").Value = "2016-03-18 16:44:30"
